# Commit: "Add files via upload"
# The bank name "Комерцијална банка АД Скопје " (with a trailing non-breaking
# space) used in column A is replaced by a clean "Комерцијална банка АД Скопје"
# (no trailing whitespace). Likewise the "Банката " (trailing space) label
# used in B21:B22 is replaced by the existing clean "Банката" label.
# Excel automatically drops the now-unused shared-string entries and appends
# the new one on save, which matches the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A23").Value = "Комерцијална банка АД Скопје"
$ws.Range("B21:B22").Value = "Банката"

# Restore the selection that was active when the workbook was saved: the
# whole of row 24 (A24:XFD1048576) with A24 as the active cell.
$ws.Range("A24:XFD1048576").Select()
